# Round the numeric result columns (B:E) for rows 2-13 to integer values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    for ($col = 2; $col -le 5; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = [Math]::Round([double]$cell.Value2, 0)
    }
}
